$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "15÷3=5, 0"
$t.Cell(1,2).Range.Text = "39÷2=19, 1"
$t.Cell(1,3).Range.Text = "57÷6=9, 3"
$t.Cell(1,4).Range.Text = "60÷4=15, 0"
$t.Cell(1,5).Range.Text = "13÷7=1, 6"

$t.Cell(5,1).Range.Text = "61÷6=10, 1"
$t.Cell(5,2).Range.Text = "51÷5=10, 1"
$t.Cell(5,3).Range.Text = "43÷9=4, 7"
$t.Cell(5,4).Range.Text = "45÷9=5, 0"
$t.Cell(5,5).Range.Text = "86÷6=14, 2"

$t.Cell(9,1).Range.Text = "87÷5=17, 2"
$t.Cell(9,2).Range.Text = "89÷3=29, 2"
$t.Cell(9,3).Range.Text = "42÷2=21, 0"
$t.Cell(9,4).Range.Text = "10÷5=2, 0"
$t.Cell(9,5).Range.Text = "21÷5=4, 1"

$t.Cell(13,1).Range.Text = "29÷6=4, 5"
$t.Cell(13,2).Range.Text = "52÷2=26, 0"
$t.Cell(13,3).Range.Text = "71÷3=23, 2"
$t.Cell(13,4).Range.Text = "26÷6=4, 2"
$t.Cell(13,5).Range.Text = "43÷8=5, 3"

$t.Cell(17,1).Range.Text = "63÷6=10, 3"
$t.Cell(17,2).Range.Text = "84÷9=9, 3"
$t.Cell(17,3).Range.Text = "38÷2=19, 0"
$t.Cell(17,4).Range.Text = "62÷3=20, 2"
$t.Cell(17,5).Range.Text = "96÷3=32, 0"
